# Update "想去人数" (want-to-go attendance counts) in column F across the
# three data sheets (展览, 演出, 全部类型). The values were refreshed from
# the live data source at build time (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 637
    3  = 6041
    6  = 1021
    7  = 407
    8  = 1383
    10 = 3135
    11 = 307
    12 = 1979
    13 = 123
    15 = 206
    16 = 88
    17 = 184
    18 = 1005
    21 = 84
    22 = 3732
    23 = 1184
    24 = 2956
    25 = 293
    26 = 2519
    27 = 4291
    29 = 937
    30 = 483
    31 = 1353
    32 = 121
    33 = 15
    34 = 45
    35 = 43
    37 = 1042
    39 = 79
    40 = 1121
    41 = 729
    42 = 641
    44 = 23
    45 = 122
    47 = 4
    48 = 325
    49 = 3622
}
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# Sheet "演出" (row -> new F value)
$ws2 = $wb.Worksheets.Item("演出")
$updates2 = @{
    2 = 23
}
foreach ($row in $updates2.Keys) {
    $ws2.Cells.Item($row, 6).Value = $updates2[$row]
}

# Sheet "全部类型" (row -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 637
    3  = 6041
    5  = 23
    6  = 407
    7  = 1383
    8  = 3135
    10 = 1979
    11 = 123
    13 = 206
    16 = 88
    17 = 184
    18 = 1005
    20 = 84
    21 = 3732
    23 = 1184
    25 = 2956
    26 = 2519
    27 = 4291
    30 = 937
    31 = 1353
    33 = 1042
    36 = 79
    37 = 1121
    39 = 729
    43 = 23
    45 = 123
    47 = 325
    48 = 3622
}
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
